$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header columns (F:J), reusing the existing header formatting -------
$ws.Cells.Item(1, 1).Copy() | Out-Null
$ws.Range("F1:J1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$newHeaders = @(
    @{Col=6;  Text="CUDA Time Numba GMC"},
    @{Col=7;  Text="CUDA Time Numba SMC"},
    @{Col=8;  Text="CUDA Time naive"},
    @{Col=9;  Text="CUDA Time Global Memory Coalescing"},
    @{Col=10; Text="CUDA Time Shared Memory Caching"}
)

foreach ($h in $newHeaders) {
    $ws.Cells.Item(1, $h.Col).Value = $h.Text
}

# --- Updated / new trial rows (A:J), rows 2-9 = Trial 1-8 --------------------
$rows = @(
    @{Row=2; Trial="Trial 1"; B=0; C=1.06029200553894;   D=0.3211402893066406; H=0.5894157886505127; I=0.2652904987335205; J=0.4007976055145264},
    @{Row=3; Trial="Trial 2"; B=0; C=1.062727451324463;  D=0.3201684951782227; H=0.5928997993469238; I=0.2629694938659668; J=0.4061641693115234},
    @{Row=4; Trial="Trial 3"; B=0; C=0.5591764450073242; D=0.1635632514953613; H=0.3131606578826904; I=0.1087098121643066; J=0.1495988368988037},
    @{Row=5; Trial="Trial 4"; B=0; C=0.5570096969604492; D=0.1581416130065918; H=0.3169970512390137; I=0.1077165603637695; J=0.1471126079559326},
    @{Row=6; Trial="Trial 5"; B=0; C=0.5584828853607178; D=0.1705427169799805; H=0.3161542415618896; I=0.1067156791687012; J=0.1486248970031738},
    @{Row=7; Trial="Trial 6"; B=0; C=0.5551903247833252; D=0.1615891456604004; H=0.3157603740692139; I=0.1097064018249512; J=0.149599552154541},
    @{Row=8; Trial="Trial 7"; B=0; C=0.5555508136749268; D=0.1625645160675049; H=0.2972052097320557; I=0.1107287406921387; J=0.1565570831298828},
    @{Row=9; Trial="Trial 8"; B=0; C=0.5511729717254639; D=0.16054368019104;   H=0.319221019744873;  I=0.105689525604248;  J=0.1506044864654541}
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $r.Trial    # A - Trial Name
    $ws.Cells.Item($row, 2).Value = $r.B        # B - Naive Time
    $ws.Cells.Item($row, 3).Value = $r.C        # C - Naive Time Numba
    $ws.Cells.Item($row, 4).Value = $r.D        # D - ikj Time Numba
    # E (CUDA Time Numba naive), F (CUDA Time Numba GMC), G (CUDA Time Numba SMC)
    # are intentionally left blank for every trial row - not measured per-trial.
    $ws.Cells.Item($row, 8).Value  = $r.H        # H - CUDA Time naive
    $ws.Cells.Item($row, 9).Value  = $r.I        # I - CUDA Time Global Memory Coalescing
    $ws.Cells.Item($row, 10).Value = $r.J        # J - CUDA Time Shared Memory Caching
}

# E2 / E3 previously held numeric "CUDA Time Numba naive" values - clear them,
# that column is no longer populated per-trial going forward.
$ws.Cells.Item(2, 5).ClearContents()
$ws.Cells.Item(3, 5).ClearContents()
